$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Word-pair dictionary data (column A / column B), rows 1-10
$data = @(
    @("brave", "history"),
    @("strong", "book"),
    @("coward", "read"),
    @("fat", "write"),
    @("ugly", "listen"),
    @("handsome", "face"),
    @("horror", "look"),
    @("fantasy", "screen"),
    @("science fiction", "computer"),
    @("comedy", "head")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

# Selection now covers the filled range
$ws.Range("A1:B10").Select()

# Reposition/resize the workbook window to match the saved view state
try {
    $win = $wb.Windows.Item(1)
    $win.Left = 768
    $win.Top = 768
    $win.Width = 14916
    $win.Height = 11808
} catch {
}
